# Update "想去人数" (want-to-go count) values on three sheets to match
# the regenerated data snapshot (gh-pages output at commit 456a3b4).

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (exhibitions) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 46
$ws1.Range("F5").Value = 1010
$ws1.Range("F9").Value = 1478
$ws1.Range("F11").Value = 1380
$ws1.Range("F13").Value = 496
$ws1.Range("F16").Value = 814
$ws1.Range("F17").Value = 250
$ws1.Range("F18").Value = 1413
$ws1.Range("F19").Value = 271
$ws1.Range("F20").Value = 66
$ws1.Range("F21").Value = 1150
$ws1.Range("F25").Value = 3569
$ws1.Range("F26").Value = 709

# --- Sheet "演出" (performances) ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F3").Value = 174
$ws2.Range("F4").Value = 35
$ws2.Range("F5").Value = 15
$ws2.Range("F7").Value = 3
$ws2.Range("F8").Value = 32
$ws2.Range("F9").Value = 22
$ws2.Range("F13").Value = 85
$ws2.Range("F14").Value = 20

# --- Sheet "全部类型" (all types combined) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 46
$ws4.Range("F7").Value = 174
$ws4.Range("F8").Value = 35
$ws4.Range("F9").Value = 15
$ws4.Range("F12").Value = 3
$ws4.Range("F13").Value = 32
$ws4.Range("F14").Value = 22
$ws4.Range("F16").Value = 1010
$ws4.Range("F20").Value = 1478
$ws4.Range("F22").Value = 1380
$ws4.Range("F24").Value = 496
$ws4.Range("F27").Value = 814
$ws4.Range("F28").Value = 250
$ws4.Range("F29").Value = 1413
$ws4.Range("F30").Value = 271
$ws4.Range("F31").Value = 66
$ws4.Range("F34").Value = 1150
$ws4.Range("F38").Value = 3569
$ws4.Range("F39").Value = 709
$ws4.Range("F42").Value = 85
$ws4.Range("F43").Value = 20
